$d = $word.ActiveDocument

$pairs = @(
  @("142÷9=15, 7", "899÷3=299, 2"),
  @("398÷7=56, 6", "162÷7=23, 1"),
  @("141÷5=28, 1", "198÷4=49, 2"),
  @("804÷8=100, 4", "432÷3=144, 0"),
  @("884÷2=442, 0", "324÷5=64, 4"),
  @("587÷6=97, 5", "842÷9=93, 5"),
  @("182÷3=60, 2", "194÷9=21, 5"),
  @("413÷9=45, 8", "308÷2=154, 0"),
  @("921÷4=230, 1", "993÷3=331, 0"),
  @("226÷5=45, 1", "988÷4=247, 0"),
  @("318÷5=63, 3", "896÷7=128, 0"),
  @("949÷8=118, 5", "964÷9=107, 1"),
  @("914÷8=114, 2", "337÷4=84, 1"),
  @("275÷3=91, 2", "340÷9=37, 7"),
  @("221÷3=73, 2", "787÷7=112, 3"),
  @("895÷8=111, 7", "602÷6=100, 2"),
  @("304÷7=43, 3", "998÷2=499, 0"),
  @("704÷8=88, 0", "128÷3=42, 2"),
  @("444÷7=63, 3", "468÷5=93, 3"),
  @("129÷6=21, 3", "867÷7=123, 6"),
  @("529÷5=105, 4", "202÷8=25, 2"),
  @("753÷9=83, 6", "728÷6=121, 2"),
  @("802÷2=401, 0", "233÷6=38, 5"),
  @("601÷6=100, 1", "661÷3=220, 1"),
  @("533÷4=133, 1", "867÷9=96, 3")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
